# Statistical significance passing tests
#
# Inserts a new "SignificanceValues" worksheet immediately before the
# existing "Lookups" worksheet. The new sheet mirrors the layout of the
# "DisplayValues" sheet (row_heading header in A1, numeric column headers
# in B1/C1, row headers 0/1 in A2/A3) but carries significance letter
# codes ("H" / "S") instead of the display percentages.

$wb = $excel.ActiveWorkbook

$displayValues = $wb.Worksheets.Item("DisplayValues")
$lookups = $wb.Worksheets.Item("Lookups")

# Duplicate DisplayValues (same column widths / formatting / header row)
# and drop the copy right before Lookups so the final sheet order becomes
# ExistingData, DisplayValues, SignificanceValues, Lookups. Excel makes
# the freshly created copy the active sheet, so grab it from there rather
# than guessing the auto-generated "DisplayValues (2)" name.
$displayValues.Copy($lookups)
$sig = $wb.ActiveSheet
$sig.Name = "SignificanceValues"

# Replace the copied 0.5 / 0.5 significance-grid placeholders with the
# real sparse significance markers.
$sig.Range("B2:C2").ClearContents()
$sig.Range("B3:C3").ClearContents()
$sig.Range("B2").Value = "H"
$sig.Range("C3").Value = "S"

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("ExistingData").Activate()
